$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date values (shifted forward by one week) and new MyForecast values
$dates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$forecasts = @(77, 79, 81, 82, 81, 81, 84, 87, 90, 91, 90, 90, 91, 94, 98, 98)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = "'" + $dates[$i]
    $ws1.Cells.Item($row, 4).Value = $forecasts[$i]
}

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# All values on the Summary sheet are stored as text, so every assignment
# is prefixed with a leading apostrophe to stop Excel from reinterpreting
# numeric- or date-looking strings as numbers/dates.
$ws2.Range("B2").Value = "'2023-02-19 to 2025-01-05"
$ws2.Range("B4").Value = "'180"
$ws2.Range("B6").Value = "'54"
$ws2.Range("B7").Value = "'39"
$ws2.Range("B8").Value = "'3075 units"
$ws2.Range("B9").Value = "'1394"
$ws2.Range("B10").Value = "'652"
$ws2.Range("B11").Value = "'319"
$ws2.Range("B12").Value = "'98"
$ws2.Range("B13").Value = "'2025-04-20"
$ws2.Range("B14").Value = "'77"
$ws2.Range("B15").Value = "'2025-01-12"
